# Daily attendance processing - 2026-01-11 14:57:36
# Reverse the order of names/emails listed in the "Recorded By" column (G)
# for each session row, leaving entries that include the protected
# "admin@admin.com" account untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val.Contains(",") -and -not $val.Contains("admin@admin.com")) {
        $parts = $val -split ",\s*"
        $n = $parts.Count
        $rev = @()
        for ($i = $n - 1; $i -ge 0; $i--) {
            $rev += $parts[$i]
        }
        $cell.Value2 = [string]::Join(", ", $rev)
    }
}
